$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 7 - State Comparison
$ws.Range("B7").Value = "State Comparison"
$ws.Range("D7").Value = "Fails as GetType () returns the super class of the type given meaning that every state is the same as each other as they all implent interface IState<T>. (fixed)"

# Fill in row 8 - GetHealthKitState
$ws.Range("B8").Value = "GetHealthKitState"
$ws.Range("D8").Value = "Failed as it never bothered to pick up the health kit. (fixed as I needed to check if the item was in the inventory to keep the loop running with the new system in place.)"

# Update selection to D8
$ws.Range("D8").Select()
